# Update the cryptos list worksheet with refreshed price/volume data
# and the two pairs of swapped rows (Frax/Aptos and Quant/NEARProtocol).
#
# NOTE: Many values in column D look like numbers (e.g. "1.001", "288.83")
# but must be stored as plain TEXT, exactly as the source data has them
# (the sheet uses inline/shared strings throughout, not numeric cells).
# To stop Excel's automatic number conversion we prefix those assignments
# with a leading apostrophe (forces text entry) and then immediately
# reset the cell style back to "Normal" so no stray formatting sticks
# around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "21.694.24"
$ws.Range("E2").Value = "  -1.50%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.533.39"
$ws.Range("E3").Value = "  -1.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - USDC
Set-TextValue "D5" "1.001"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6 - BNB
Set-TextValue "D6" "288.83"
$ws.Range("E6").Value = "  +0.58%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3917"
$ws.Range("E7").Value = "  +2.94%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3159"
$ws.Range("E8").Value = "  -2.74%  "

# Row 9 - OKB
Set-TextValue "D9" "42.28"
$ws.Range("E9").Value = "  +1.73%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.07170"
$ws.Range("E10").Value = "  -2.13%  "

# Row 11 - Polygon
Set-TextValue "D11" "1.045"
$ws.Range("E11").Value = "  -6.82%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  +0.10%  "

# Row 13 - Polkadot
Set-TextValue "D13" "5.630"
$ws.Range("E13").Value = "  -1.67%  "

# Row 14 - Solana
Set-TextValue "D14" "18.52"
$ws.Range("E14").Value = "  -4.56%  "

# Row 15 - Chainlink
Set-TextValue "D15" "6.591"
$ws.Range("E15").Value = "  -3.26%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.538.79"
$ws.Range("E16").Value = "  -1.18%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.00001097"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18 - TRON
Set-TextValue "D18" "0.06584"
$ws.Range("E18").Value = "  -0.55%  "

# Row 19 - Litecoin
Set-TextValue "D19" "83.00"
$ws.Range("E19").Value = "  -2.49%  "

# Row 20 - Dai
Set-TextValue "D20" "1.001"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.112"
$ws.Range("E21").Value = "  -4.88%  "

# Row 22 - Avalanche
Set-TextValue "D22" "15.40"
$ws.Range("E22").Value = "  -3.59%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -5.74%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.375"
$ws.Range("E24").Value = "  +3.27%  "

# Row 25 - WrappedBTC
$ws.Range("D25").Value = "21.703.97"
$ws.Range("E25").Value = "  -1.51%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.352"
$ws.Range("E26").Value = "  -6.92%  "

# Row 27 - Monero
Set-TextValue "D27" "147.46"
$ws.Range("E27").Value = "  -0.80%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "18.32"
$ws.Range("E28").Value = "  -2.86%  "

# Row 29 - HuobiToken
Set-TextValue "D29" "4.843"
$ws.Range("E29").Value = "  -0.11%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.713.73"
$ws.Range("E30").Value = "  -1.06%  "

# Row 31 - BitcoinCash
Set-TextValue "D31" "117.17"
$ws.Range("E31").Value = "  -2.80%  "

# Row 32 - Filecoin
Set-TextValue "D32" "5.878"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33 - ImmutableX
Set-TextValue "D33" "0.9584"
$ws.Range("E33").Value = "  -13.08%  "

# Row 34 - Stellar
Set-TextValue "D34" "0.08140"
$ws.Range("E34").Value = "  -0.15%  "

# Row 35 - FraxShare
Set-TextValue "D35" "8.637"
$ws.Range("E35").Value = "  -7.09%  "

# Row 36 - Hedera
Set-TextValue "D36" "0.06065"
$ws.Range("E36").Value = "  -2.27%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue "D37" "5.117"
$ws.Range("E37").Value = "  -2.52%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.02202"
$ws.Range("E38").Value = "  -3.78%  "

# Row 39 - WEMIXTOKEN
Set-TextValue "D39" "1.452"
$ws.Range("E39").Value = "  -11.73%  "

# Row 40 - Algorand
Set-TextValue "D40" "0.2025"
$ws.Range("E40").Value = "  -3.79%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "1.179"
$ws.Range("E41").Value = "  -3.54%  "

# Row 42/43 - Frax and Aptos swap places (Aptos now row 42, Frax now row 43)
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D42" "10.88"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D43" "1.000"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44 - TheSandbox
Set-TextValue "D44" "0.5755"
$ws.Range("E44").Value = "  -3.11%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.00"
$ws.Range("E45").Value = "  -3.84%  "

# Row 46 - PancakeSwap
Set-TextValue "D46" "3.730"

# Row 47 - Decentraland
Set-TextValue "D47" "0.5486"
$ws.Range("E47").Value = "  -4.44%  "

# Row 48 - EOS
Set-TextValue "D48" "1.152"
$ws.Range("E48").Value = "  -0.36%  "

# Row 49/50 - Quant and NEARProtocol swap places (NEARProtocol now row 49, Quant now row 50)
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "1.865"
$ws.Range("E49").Value = "  -3.56%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "115.72"
$ws.Range("E50").Value = "  -3.28%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.06688"
